$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new test data row (row 4)
$ws.Range("B4").Value = "Pune"
$ws.Range("C4").Value = "Hyderabad"
$ws.Range("D4").Value = "07-Jun-2021"

# Best-fit the new column C like the other text columns
$ws.Columns.Item(3).AutoFit() | Out-Null

# Match the saved cursor/selection position
$ws.Range("H7").Select() | Out-Null
